$wb = $excel.ActiveWorkbook

# Sheet "展览" (Exhibitions)
$ws1 = $wb.Worksheets.Item("展览")
$ws1.Range("F2").Value = 77
$ws1.Range("F8").Value = 4387
$ws1.Range("F10").Value = 4955
$ws1.Range("F11").Value = 552
$ws1.Range("F12").Value = 1248
$ws1.Range("F13").Value = 85

# Sheet "全部类型" (All types)
$ws4 = $wb.Worksheets.Item("全部类型")
$ws4.Range("F2").Value = 77
$ws4.Range("F9").Value = 4387
$ws4.Range("F11").Value = 4955
$ws4.Range("F12").Value = 552
$ws4.Range("F13").Value = 1248
$ws4.Range("F14").Value = 85
